# BOT; UPDATE DATA
# Appends the 2020-05-05 (Excel serial 43956) daily COVID-19 figures to the
# three data sheets ("all", "kobe", "other"), pushing each sheet's trailing
# footnote row(s) down by one row, same as the upstream data-refresh bot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert new row 28 (date 2020-05-05), footnote moves 28->29
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("all")
$ws.Activate()

# Push the footnote (currently B28) down to B29, carrying its formatting.
$ws.Range("B28").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Cells.Item(29, 2).Value2 = $ws.Cells.Item(28, 2).Value2

# Build the new data row 28 off row 27's formatting, then fill the values.
$ws.Range("A27:H27").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)
$ws.Cells.Item(28, 1).Value2 = 43956
$ws.Cells.Item(28, 2).Value2 = 269
$ws.Cells.Item(28, 3).Value2 = 263
$ws.Cells.Item(28, 4).Value2 = 121
$ws.Cells.Item(28, 5).Value2 = 112
$ws.Cells.Item(28, 6).Value2 = 9
$ws.Cells.Item(28, 7).Value2 = 7
$ws.Cells.Item(28, 8).Value2 = 135

$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 2

# ---------------------------------------------------------------------
# Sheet "kobe": correct D82, insert new row 83, blank/footnote rows
# shift 83->84 and 84->85
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("kobe")
$ws.Activate()

# Correction to the already-published 2020-05-04 row.
$ws.Cells.Item(82, 4).Value2 = 2

# Push the footnote (currently B84) down to B85.
$ws.Range("B84").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$ws.Cells.Item(85, 2).Value2 = $ws.Cells.Item(84, 2).Value2
$ws.Cells.Item(84, 2).Clear()

# Push the (then-empty) placeholder row, currently A83, down to row 84.
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)

# Build the new data row 83 off row 82's formatting, then fill the values.
$ws.Range("A82:J82").Copy()
$ws.Range("A83:J83").PasteSpecial(-4122)
$ws.Cells.Item(83, 1).Value2 = 43956
$ws.Cells.Item(83, 2).ClearContents()
$ws.Cells.Item(83, 3).Value2 = 2242
$ws.Cells.Item(83, 4).Value2 = 1
$ws.Cells.Item(83, 5).Value2 = 269
$ws.Cells.Item(83, 6).Value2 = 116
$ws.Cells.Item(83, 7).Value2 = 108
$ws.Cells.Item(83, 8).Value2 = 8
$ws.Cells.Item(83, 9).Value2 = 7
$ws.Cells.Item(83, 10).Value2 = 128

$ws.Range("G84").Select()
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.ScrollColumn = 2

# ---------------------------------------------------------------------
# Sheet "other": insert new row 58, blank/footnote rows shift 58->59
# and 59->60
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("other")
$ws.Activate()

# Push the footnote (currently B59) down to B60.
$ws.Range("B59").Copy()
$ws.Range("B60").PasteSpecial(-4122)
$ws.Cells.Item(60, 2).Value2 = $ws.Cells.Item(59, 2).Value2
$ws.Cells.Item(59, 2).Clear()

# Push the (then-blank) placeholder row, currently A58:I58, down to row 59.
$ws.Range("A58:I58").Copy()
$ws.Range("A59:I59").PasteSpecial(-4122)

# Build the new data row 58 off row 57's formatting (A:H; I58 already
# carried its blank formatting down with the rest of the row), then fill.
$ws.Range("A57:H57").Copy()
$ws.Range("A58:H58").PasteSpecial(-4122)
$ws.Cells.Item(58, 1).Value2 = 43956
$ws.Cells.Item(58, 2).Value2 = 0
$ws.Cells.Item(58, 3).Value2 = 12
$ws.Cells.Item(58, 4).Value2 = 5
$ws.Cells.Item(58, 5).Value2 = 4
$ws.Cells.Item(58, 6).Value2 = 1
$ws.Cells.Item(58, 7).Value2 = 0
$ws.Cells.Item(58, 8).Value2 = 7

$ws.Range("G59").Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 4

# ---------------------------------------------------------------------
# Leave the workbook on the sheet that was active originally.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("all").Activate()
